$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Förändrad" (Changed) date column (C2:C13) from 45170 to 45174
# (equivalent to updating the date from 2023-09-01 to 2023-09-05)
$ws.Range("C2:C13").Value = 45174
